$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-29: row number followed by values for columns D, J, K, L, M, N, O, P, Q
$data = @(
    @(2, 44302, 300, 900, 1000, 950, '$/atado 1,5 a 2 kilos', 'Región de Arica y Parinacota', 475, 2),
    @(3, 44181, 200, 1000, 1200, 1100, '$/atado', 'Región de Arica y Parinacota', 1100, 1),
    @(4, 44253, 250, 1800, 2000, 1900, '$/atado 1,5 a 2 kilos', 'Región de Arica y Parinacota', 950, 2),
    @(5, 44365, 200, 1800, 2000, 1900, '$/atado 1,5 a 2 kilos', 'Región de Arica y Parinacota', 950, 2),
    @(6, 44544, 250, 900, 1000, 950, '$/atado 1,5 a 2 kilos', 'Región de Arica y Parinacota', 475, 2),
    @(7, 44438, 300, 950, 1000, 975, '$/atado 1,5 a 2 kilos', 'Región de Arica y Parinacota', 488, 2),
    @(8, 44540, 300, 900, 1000, 950, '$/atado 1,5 a 2 kilos', 'Región de Arica y Parinacota', 475, 2),
    @(9, 44363, 250, 2500, 2800, 2650, '$/atado 1,5 a 2 kilos', 'Región de Arica y Parinacota', 1325, 2),
    @(10, 44447, 300, 900, 1000, 950, '$/atado 1,5 a 2 kilos', 'Región de Arica y Parinacota', 475, 2),
    @(11, 44403, 250, 1800, 2000, 1900, '$/atado 1,5 a 2 kilos', 'Región de Arica y Parinacota', 950, 2),
    @(12, 44172, 200, 1300, 1500, 1400, '$/atado 1,5 a 2 kilos', 'Región de Arica y Parinacota', 700, 2),
    @(13, 44229, 250, 1800, 2000, 1900, '$/atado 1,5 a 2 kilos', 'Región de Arica y Parinacota', 950, 2),
    @(14, 44390, 250, 2400, 2500, 2450, '$/atado 1,5 a 2 kilos', 'Región de Arica y Parinacota', 1225, 2),
    @(15, 44427, 250, 1300, 1500, 1400, '$/atado 1,5 a 2 kilos', 'Región de Arica y Parinacota', 700, 2),
    @(16, 44266, 300, 1700, 1800, 1750, '$/atado 1,5 a 2 kilos', 'Región de Arica y Parinacota', 875, 2),
    @(17, 44392, 250, 1800, 2000, 1900, '$/atado 1,5 a 2 kilos', 'Región de Arica y Parinacota', 950, 2),
    @(18, 44572, 300, 1400, 1500, 1450, '$/atado 1,5 a 2 kilos', 'Región de Arica y Parinacota', 725, 2),
    @(19, 44616, 270, 1300, 1500, 1400, '$/atado 1,5 a 2 kilos', 'Región de Arica y Parinacota', 700, 2),
    @(20, 44257, 500, 1400, 1500, 1450, '$/atado 1,5 a 2 kilos', 'Región de Arica y Parinacota', 725, 2),
    @(21, 44385, 300, 2400, 2500, 2450, '$/atado 1,5 a 2 kilos', 'Región de Arica y Parinacota', 1225, 2),
    @(22, 44468, 300, 900, 1000, 950, '$/atado 1,5 a 2 kilos', 'Región de Arica y Parinacota', 475, 2),
    @(23, 44161, 270, 900, 1000, 950, '$/atado 1,5 a 2 kilos', 'Región de Arica y Parinacota', 475, 2),
    @(24, 44601, 270, 2200, 2500, 2350, '$/atado 1,5 a 2 kilos', 'Región de Arica y Parinacota', 1175, 2),
    @(25, 44525, 300, 1400, 1500, 1450, '$/atado 1,5 a 2 kilos', 'Región de Arica y Parinacota', 725, 2),
    @(26, 44291, 250, 1800, 2000, 1900, '$/atado 1,5 a 2 kilos', 'Región de Arica y Parinacota', 950, 2),
    @(27, 44435, 300, 900, 1000, 950, '$/atado 1,5 a 2 kilos', 'Región de Arica y Parinacota', 475, 2),
    @(28, 44202, 250, 1800, 2000, 1900, '$/atado 1,5 a 2 kilos', 'Región de Arica y Parinacota', 950, 2),
    @(29, 44243, 250, 1200, 1300, 1250, '$/atado 1,5 a 2 kilos', 'Región de Arica y Parinacota', 625, 2),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("D$r").Value = $row[1]
    $ws.Range("J$r").Value = $row[2]
    $ws.Range("K$r").Value = $row[3]
    $ws.Range("L$r").Value = $row[4]
    $ws.Range("M$r").Value = $row[5]
    $ws.Range("N$r").Value = $row[6]
    $ws.Range("O$r").Value = $row[7]
    $ws.Range("P$r").Value = $row[8]
    $ws.Range("Q$r").Value = $row[9]
}

Write-Host "Applied weekly fruit/vegetable price update to rows 2-29."
